$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Strings that already exist in the sheet and are simply being moved to a new
# row (Docentes responsaveis / Programa resumido / Programa / Avaliacao block
# all shift down by two rows to make room for the new TG-specific content).
# ---------------------------------------------------------------------------
$sDocente1            = "5840730 - Antonio Jefferson da Silva Machado"
$sDocente2            = "1176388 - Luiz Tadeu Fernandes Eleno"
$sProgResumido        = "Programa resumido:"
$sShortSyllabus       = "Short syllabus:"
$sPreparation         = "Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology."
$sPrograma            = "Programa:"
$sSyllabus            = "Syllabus:"
$sStudentShouldLook   = "The student should look for a professor or professional with training in the area of engineering or related areas, for the elaboration of a project proposal containing motivation and objectives, theoretical foundation and execution schedule. The project itself will be developed and defended in the course Undergraduate Work II."
$sAvaliacao           = "Avaliação:"
$sMetodo              = "Método:"
$sAlunoApresentar     = "O aluno deve apresentar a proposta de trabalho à uma banca formada pelo responsável pela disciplina e professores ou profissionais da área."
$sCriterio            = "Critério:"
$sAvaliacaoAtribuicao = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."
$sNormaRecuperacao    = "Norma de recuperação:"
$sCriterioBanca       = "A critério da banca de avaliação poderá ser estabelecido um prazo para readequação e reapresentação do plano de trabalho."
$sRequisitos          = "Requisitos:"
$sLOM3238             = "LOM3238 -  Projeto Integrado  (Requisito)`n"
$sBibliografia        = "Bibliografia:"

# ---------------------------------------------------------------------------
# Brand new strings introduced by this revision.
# ---------------------------------------------------------------------------
$sObjetivosPt   = "O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico."
$sElaboracao    = "Elaboração, com a orientação de um professor supervisor, de uma proposta de projeto em tema ligado à área de ciência e tecnologia."
$sAlunoProcurar = "O aluno deve procurar um professor ou profissional com formação na área de engenharia ou áreas correlatas, para a elaboração de uma proposta de projeto contendo motivação e objetivos, fundamentação teórica e cronograma de execução. O projeto propriamente dito será desenvolvido e defendido na disciplina Trabalho de Graduação II."
$sASerDefinida  = "A ser definida no plano de trabalho."

# Row 10: Objetivos (PT) gains the real text (was a stray copy of the
# "Docentes responsaveis" value).
$ws.Range("B10").Value = $sObjetivosPt
$ws.Range("C10").Value = $sObjetivosPt

# Row 13: now holds just the first docente (no row label, default height).
$ws.Range("A13").Clear()
$ws.Range("B13").Value = $sDocente1
$ws.Range("C13").Value = $sDocente1
$ws.Rows.Item(13).AutoFit()

# Row 14: now holds just the second docente (no row label, default height).
$ws.Range("A14").Clear()
$ws.Range("B14").Value = $sDocente2
$ws.Range("C14").Value = $sDocente2
$ws.Rows.Item(14).AutoFit()

# Row 15: Programa resumido / Elaboracao text.
$ws.Range("A15").Value = $sProgResumido
$ws.Range("B15").Value = $sElaboracao
$ws.Range("C15").Value = $sElaboracao
$ws.Rows.Item(15).RowHeight = 60

# Row 16: Short syllabus / Preparation text.
$ws.Range("A16").Value = $sShortSyllabus
$ws.Range("B16").Value = $sPreparation
$ws.Range("C16").Value = $sPreparation
$ws.Rows.Item(16).RowHeight = 60

# Row 17: Programa / new "aluno deve procurar..." text. Column B did not
# exist on this row before, so pull its number format (wrap, non-bold) from
# an existing column-B cell before writing the value.
$ws.Range("A17").Value = $sPrograma
$ws.Range("B18").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = $sAlunoProcurar
$ws.Range("C17").Value = $sAlunoProcurar
$ws.Rows.Item(17).RowHeight = 120

# Row 18: Syllabus / student-should-look text.
$ws.Range("A18").Value = $sSyllabus
$ws.Range("B18").Value = $sStudentShouldLook
$ws.Range("C18").Value = $sStudentShouldLook
$ws.Rows.Item(18).RowHeight = 120

# Row 19: Avaliacao label only again (no height, no B/C).
$ws.Range("A19").Value = $sAvaliacao
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()
$ws.Rows.Item(19).AutoFit()

# Row 20: Metodo / aluno apresentar text.
$ws.Range("A20").Value = $sMetodo
$ws.Range("B20").Value = $sAlunoApresentar
$ws.Range("C20").Value = $sAlunoApresentar
$ws.Rows.Item(20).RowHeight = 60

# Row 21: Criterio / avaliacao-atribuicao text.
$ws.Range("A21").Value = $sCriterio
$ws.Range("B21").Value = $sAvaliacaoAtribuicao
$ws.Range("C21").Value = $sAvaliacaoAtribuicao
$ws.Rows.Item(21).RowHeight = 60

# Row 22: Norma de recuperacao / criterio-banca text. Column B did not exist
# on this row before either, so copy formatting the same way as row 17.
$ws.Range("A22").Value = $sNormaRecuperacao
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = $sCriterioBanca
$ws.Range("C22").Value = $sCriterioBanca
$ws.Rows.Item(22).RowHeight = 60

# Row 23: Bibliografia / "a ser definida" text.
$ws.Range("A23").Value = $sBibliografia
$ws.Range("B23").Value = $sASerDefinida
$ws.Range("C23").Value = $sASerDefinida
$ws.Rows.Item(23).RowHeight = 120

# Row 24: Requisitos label, new row (label only).
$ws.Range("A24").Value = $sRequisitos

# Row 25: the requirement text, new row (B/C only). Copy formatting from
# B23 before writing the value.
$ws.Range("B23").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = $sLOM3238
$ws.Range("C25").Value = $sLOM3238
$ws.Rows.Item(25).RowHeight = 30
